$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 5) mirroring the structure of existing rows.
# Copy formatting from the row above (keeps the shared date-format style index)
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 42647.680694444447

$ws.Range("B5").Value = $false

$ws.Range("C5").Value = 9300.0400000000009
$ws.Range("D5").Value = 9646.34
$ws.Range("E5").Value = 18.12
$ws.Range("F5").Value = 18.77

$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = $true

$ws.Range("H5").Value = 3.59
$ws.Range("I5").Value = $true
